# Applies the Simplified-Chinese translation edits described by the diff.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND:" $old
    }
}

# Title
Replace-Text "提醒 ROW 客戶的電子郵件" "向 ROW 客户发送提醒电子邮件"

# "Subject line" label
Replace-Text "主題：" "Subject line:"

# Subject text
Replace-Text "我們將於 9 月 29 日移除 Tether Omni (USDT)" "将于 9 月 29 日移除 Tether Omni (USDT)"

# "Body" label
Replace-Text "本文：" "Body:"

# Heading "Saying goodbye to Tether Omni"
Replace-Text "向 Tether Omni 道別" "向 Tether Omni 说再见"

# Main paragraph about removal
Replace-Text "我們將停止在 Deriv 上提供 Tether Omni (USDT) 作為帳戶貨幣，自 2023/09/29 (00:00 GMT) 起生效。 這是因為 Tether 停止支持 USDT 轉帳的 Omni。" "自 2023 年 9 月 29 日格林威治标准时间 00:00 起，Deriv 将停止提供Tether Omni (USDT) 作为账户货币。 这是因为Tether已停止支持USDT的Omni转账。"

# "What do I need to do?" heading run
Replace-Text "我需要做什麼？" "我需要做什么？"

# First sentence run (no longer ends with trailing space)
Replace-Text "如果您的 USDT 帳戶有餘額 " "如果您在USDT账户"

# Bold account id placeholder
Replace-Text "[帳戶 ID]" "[账户ID]"

# Remaining sentence run
Replace-Text "，請在上述日期之前提取您的餘額。 如果您有未平倉的頭寸，請在提取餘額之前先關閉它們。" "中有余额，请在上述日期之前提取余额。 如果您有未平仓头寸，请在提取余额之前先关闭它们。"

# Button text
Replace-Text "檢查我的帳戶" "查看我的账户"

# "Important" paragraph main text
Replace-Text "您的 USDT 帳戶將於 2023/09/29 00:00 GMT 關閉。 任何未平倉頭寸將自動關閉，帳戶餘額將在上述日期後轉移到您最後一個有效的帳戶中" "您的USDT账户将在2023年9月29日00:00 GMT关闭。 任何未平仓头寸将在提到的日期后自动关闭，账户余额将转移到您最后活跃的账户"

# Trailing sentence after comment references
Replace-Text "在此過程中將適用標準匯率和費用。" "在此过程中将适用标准汇率和费用。"

# "Contact us" line
Replace-Text "如有任何疑問，請聯繫我們：" "如有任何疑问，请联系我们："

# "Live chat" hyperlink text
Replace-Text "即時聊天" "实时聊天"

# Comment 0, first paragraph
Replace-Text "@azita@regentmarkets.com，BE 不能保證他們能在那時準備好腳本，" "@azita@regentmarkets.com，BE不能保证他们可以在那时完成脚本，"

# Comment 0, second paragraph
Replace-Text "我們可以說成「在提到的日期之後」嗎？" "我们可以提到 ""提到的日期后"" 吗？"

# Comment 1
Replace-Text "您是指轉帳將在提到的日期之後進行嗎？" "您是指在提到的日期之后完成转账吗？"

# Comment 2
Replace-Text "是的..看起來我們無法確認日期" "是的..似乎我们无法确认日期"

Write-Host "Done"
